# Apply updated values to the active worksheet to match the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.974
$ws.Range("D4").Value = -8.25

$ws.Range("D5").Value = -8.581999999999999

$ws.Range("C6").Value = -12.445

$ws.Range("C7").Value = -12.673

$ws.Range("D8").Value = -8.196

$ws.Range("C16").Value = -12.072
$ws.Range("D16").Value = -8.574000000000002

$ws.Range("C20").Value = -13.041

$ws.Range("D22").Value = -8.177000000000001
